# Add two new rows (Rotary Encoder PinA / PinB) to the "Frameside" pinout sheet,
# just above the existing "High Power Solenoid" row, shifting the rows below down.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Frameside")

# Insert two blank rows at 14:15 (existing rows 14+ shift down to 16+).
$ws.Range("A14:G15").Insert(-4121)

# The newly inserted rows don't inherit formatting automatically; copy the
# format from row 16 (the row that used to be row 14, "High Power Solenoid",
# which still carries the correct style for this table) onto the two blanks.
$ws.Range("A16:G16").Copy()
$ws.Range("A14:G15").PasteSpecial(-4122)

# Row 14: Rotary Encoder PinA
$ws.Range("A14").Value = "Rotary Encoder PinA"
$ws.Range("B14").Value = "Digital"
$ws.Range("C14").Value = "Input"
$ws.Range("D14").Value = 20
$ws.Range("E14").Value = "No"
$ws.Range("F14").Value = "Yes(3)"
$ws.Range("G14").Value = "Interrupt Pin A on the rotary encoder."

# Row 15: Rotary Encoder PinB
$ws.Range("A15").Value = "Rotary Encoder PinB"
$ws.Range("B15").Value = "Digital"
$ws.Range("C15").Value = "Input"
$ws.Range("D15").Value = 21
$ws.Range("E15").Value = "No"
$ws.Range("F15").Value = "Yes(2)"
$ws.Range("G15").Value = "Interrupt Pin B on the rotary encoder."

# Both new rows use a shorter row height than the wrapped-text rows around them.
$ws.Rows.Item(14).RowHeight = 14.9
$ws.Rows.Item(15).RowHeight = 14.9

# Column A got wider to fit the new, longer signal names.
$ws.Columns.Item(1).ColumnWidth = 18

# Selection moved onto the newly added Pin B description cell.
$ws.Range("G15").Select() | Out-Null

Write-Host "Rotary encoder pins added to Frameside sheet"
